# "unir_datos" edit: merge/reshape the raw object/time/distance table into
# a Tiempo / distancia / acceleration / velocidad series table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the leading "objeto" column (A) and the trailing "objeto_distancia"
# column (originally F, now E after the shift) so only the four measurement
# columns remain.
$ws.Columns("A:A").Delete()
$ws.Columns("E:E").Delete()

# Rewrite the header row.
$ws.Range("A1").Value = "Tiempo"
$ws.Range("B1").Value = "distancia"
$ws.Range("C1").Value = "acceleration"
$ws.Range("D1").Value = "velocidad"

# Rebuild the data rows with the merged/derived series.
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = -20
$ws.Range("D2").Value = 25

$ws.Range("A3").Value = 2
$ws.Range("C3").Value = -10
$ws.Range("D3").Value = 20

$ws.Range("A4").Value = 3
$ws.Range("C4").Value = -6.666666666666667
$ws.Range("D4").Value = 15

$ws.Range("A5").Value = 4
$ws.Range("C5").Value = -5
$ws.Range("D5").Value = 10

# Column B ("distancia") is left blank for every data row in the merged
# result.
$ws.Range("B2:B5").Value = ""
